$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 3; $r -le 83; $r++) {
    $ws.Cells.Item($r, 2).Value = 32
}
